# Add season-record columns (Wins / Losses / Ties) to the COL_2008 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, centered, thin border - same style
# used by all other header cells in row 1) onto the three new header cells.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-51) gets the team's season record: 74 wins, 88
# losses, 0 ties.
$ws.Range("AD2:AD51").Value = 74
$ws.Range("AE2:AE51").Value = 88
$ws.Range("AF2:AF51").Value = 0
